$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H594").Value = 0.38011850620757043
$ws.Range("I594").Value = 0.32156712638572876
$ws.Range("J594").Value = 0.9123407056321989
$ws.Range("K594").Value = 0.9706896551724138
$ws.Range("L594").Value = 0.6902010050251256

$ws.Range("H595").Value = 0.3881853480056075
$ws.Range("I595").Value = 0.3270220117558246
$ws.Range("J595").Value = 0.9104111829870422
$ws.Range("K595").Value = 0.9716666666666667
$ws.Range("L595").Value = 0.6910691823899371

$ws.Range("H597").Value = 0.49903290145720575
$ws.Range("I597").Value = 0.4126770350680919
$ws.Range("J597").Value = 0.8844243231580867
$ws.Range("K597").Value = 0.9351851851851852
$ws.Range("L597").Value = 0.6958841463414634

$ws.Range("H598").Value = 0.5097779619182066
$ws.Range("I598").Value = 0.42104244144356545
$ws.Range("J598").Value = 0.8807567888857353
$ws.Range("K598").Value = 0.9321824907521579
$ws.Range("L598").Value = 0.6979412327108165

$ws.Range("H600").Value = 0.5213826062484477
$ws.Range("I600").Value = 0.41821880756241503
$ws.Range("J600").Value = 0.5861423656810953
$ws.Range("K600").Value = 0.6986374565856265
$ws.Range("L600").Value = 0.7742415928509838

$ws.Range("H601").Value = 0.48637309446461324
$ws.Range("I601").Value = 0.38869219461668386
$ws.Range("J601").Value = 0.540372168283071
$ws.Range("K601").Value = 0.6471241394625805
$ws.Range("L601").Value = 0.7848628823872634

$ws.Range("H603").Value = 0.3823497878761381
$ws.Range("I603").Value = 0.32309551718907115
$ws.Range("J603").Value = 0.9123407056321989
$ws.Range("K603").Value = 0.9706896551724138
$ws.Range("L603").Value = 0.7263819095477387

$ws.Range("H604").Value = 0.3904948922600501
$ws.Range("I604").Value = 0.32857834549222686
$ws.Range("J604").Value = 0.9104111829870422
$ws.Range("K604").Value = 0.9730639730639731
$ws.Range("L604").Value = 0.7270841192003018

$ws.Range("H606").Value = 0.5021747692515613
$ws.Range("I606").Value = 0.41466692112629105
$ws.Range("J606").Value = 0.8844243231580867
$ws.Range("K606").Value = 0.9351851851851852
$ws.Range("L606").Value = 0.7324695121951219

$ws.Range("H607").Value = 0.5130073192020291
$ws.Range("I607").Value = 0.4230734575625273
$ws.Range("J607").Value = 0.8807567888857353
$ws.Range("K607").Value = 0.9355246523388117
$ws.Range("L607").Value = 0.7339741549430263

$ws.Range("H609").Value = 0.525131509760614
$ws.Range("I609").Value = 0.42026752851169413
$ws.Range("J609").Value = 0.5861423656810953
$ws.Range("K609").Value = 0.71001221001221
$ws.Range("L609").Value = 0.8039171203871749

$ws.Range("H610").Value = 0.48990977323459933
$ws.Range("I610").Value = 0.3906029489080301
$ws.Range("J610").Value = 0.540372168283071
$ws.Range("K610").Value = 0.6867369589345172
$ws.Range("L610").Value = 0.8104838709677419

$ws.Range("H612").Value = 0.38455415334605153
$ws.Range("I612").Value = 0.3246303998019212
$ws.Range("J612").Value = 0.9123407056321989
$ws.Range("K612").Value = 0.9705882352941176
$ws.Range("L612").Value = 0.7518527823137797

$ws.Range("H613").Value = 0.3927782330341543
$ws.Range("I613").Value = 0.3301422006966366
$ws.Range("J613").Value = 0.9104111829870422
$ws.Range("K613").Value = 0.9713322091062394
$ws.Range("L613").Value = 0.752561765260577

$ws.Range("H615").Value = 0.5052850907873817
$ws.Range("I615").Value = 0.4166629113791644
$ws.Range("J615").Value = 0.8844243231580867
$ws.Range("K615").Value = 0.9463806970509383
$ws.Range("L615").Value = 0.7584105623968516

$ws.Range("H616").Value = 0.5162045949918995
$ws.Range("I616").Value = 0.4251111185434582
$ws.Range("J616").Value = 0.8807567888857353
$ws.Range("K616").Value = 0.9427083333333334
$ws.Range("L616").Value = 0.7592168827866769

$ws.Range("H618").Value = 0.528934320729364
$ws.Range("I618").Value = 0.4223176577092402
$ws.Range("J618").Value = 0.5861423656810953
$ws.Range("K618").Value = 0.7399527186761229
$ws.Range("L618").Value = 0.8285693182657508

$ws.Range("H619").Value = 0.49350130809930975
$ws.Range("I619").Value = 0.3925141653334392
$ws.Range("J619").Value = 0.540372168283071
$ws.Range("K619").Value = 0.7238154613466334
$ws.Range("L619").Value = 0.8352392416491122

$ws.Range("H621").Value = 0.3929579321024846
$ws.Range("I621").Value = 0.3276180797338138
$ws.Range("J621").Value = 0.9123407056321989
$ws.Range("K621").Value = 0.9633401221995926
$ws.Range("L621").Value = 0.7953026422637266

$ws.Range("H622").Value = 0.39081888883338844
$ws.Range("I622").Value = 0.3255608810195392
$ws.Range("J622").Value = 0.9248920645898275
$ws.Range("K622").Value = 0.9638095238095238
$ws.Range("L622").Value = 0.7969327073552426

$ws.Range("H624").Value = 0.5167676068138207
$ws.Range("I624").Value = 0.4205829469392376
$ws.Range("J624").Value = 0.8844243231580867
$ws.Range("K624").Value = 0.9451303155006858
$ws.Range("L624").Value = 0.8059095808762919

$ws.Range("H625").Value = 0.5179632945197977
$ws.Range("I625").Value = 0.4223774590907704
$ws.Range("J625").Value = 0.8819629755320219
$ws.Range("K625").Value = 0.9407894736842105
$ws.Range("L625").Value = 0.8071791613722998

$ws.Range("H627").Value = 0.5418602944533498
$ws.Range("I627").Value = 0.4263936675085872
$ws.Range("J627").Value = 0.5861423656810953
$ws.Range("K627").Value = 0.7674810446503791
$ws.Range("L627").Value = 0.8635140874982302

$ws.Range("H628").Value = 0.5438910092497575
$ws.Range("I628").Value = 0.4281347906443249
$ws.Range("J628").Value = 0.5878746147889197
$ws.Range("K628").Value = 0.7302816901408451
$ws.Range("L628").Value = 0.8773060029282577

$ws.Range("H630").Value = 0.43445002942996563
$ws.Range("I630").Value = 0.30741873339957565
$ws.Range("J630").Value = 0.828018210055704
$ws.Range("K630").Value = 0.8909090909090909
$ws.Range("L630").Value = 0.8842406876790831

$ws.Range("H631").Value = 0.446704531147796
$ws.Range("I631").Value = 0.3125357315394587
$ws.Range("J631").Value = 0.8308291261911909
$ws.Range("K631").Value = 0.8894472361809045
$ws.Range("L631").Value = 0.8892520669102095

$ws.Range("H633").Value = 0.5117645949022204
$ws.Range("I633").Value = 0.3522441977400258
$ws.Range("J633").Value = 0.733045202648451
$ws.Range("K633").Value = 0.8389513108614233
$ws.Range("L633").Value = 0.8969413598285603

$ws.Range("H634").Value = 0.5168421340123284
$ws.Range("I634").Value = 0.35383897393016095
$ws.Range("J634").Value = 0.7290896809931678
$ws.Range("K634").Value = 0.8303249097472925
$ws.Range("L634").Value = 0.8979113800507516

$ws.Range("H636").Value = 0.43160108295277955
$ws.Range("I636").Value = 0.29063396977968703
$ws.Range("J636").Value = 0.5280575399523143
$ws.Range("K636").Value = 0.5951972555746141
$ws.Range("L636").Value = 0.9157151754203862

$ws.Range("H637").Value = 0.4134759627671508
$ws.Range("I637").Value = 0.2769494821753848
$ws.Range("J637").Value = 0.5081153680062789
$ws.Range("K637").Value = 0.5681470137825421
$ws.Range("L637").Value = 0.9195281230250685

$ws.Range("H639").Value = 0.43785391520307454
$ws.Range("I639").Value = 0.30889502730964874
$ws.Range("J639").Value = 0.828018210055704
$ws.Range("K639").Value = 0.8848484848484849
$ws.Range("L639").Value = 0.897803247373448

$ws.Range("H640").Value = 0.4503908164834533
$ws.Range("I640").Value = 0.3140403711042752
$ws.Range("J640").Value = 0.8308291261911909
$ws.Range("K640").Value = 0.8882978723404256
$ws.Range("L640").Value = 0.9013814274750576

$ws.Range("H642").Value = 0.5162265243855738
$ws.Range("I642").Value = 0.3539531736819643
$ws.Range("J642").Value = 0.733045202648451
$ws.Range("K642").Value = 0.8345864661654135
$ws.Range("L642").Value = 0.9105960264900662

$ws.Range("H643").Value = 0.5214534842747636
$ws.Range("I643").Value = 0.355558188471215
$ws.Range("J643").Value = 0.7290896809931678
$ws.Range("K643").Value = 0.85
$ws.Range("L643").Value = 0.9105058365758755

$ws.Range("H645").Value = 0.4356066086437315
$ws.Range("I645").Value = 0.2920400533087781
$ws.Range("J645").Value = 0.5280575399523143
$ws.Range("K645").Value = 0.6203007518796992
$ws.Range("L645").Value = 0.9278964667214462

$ws.Range("H646").Value = 0.4173819763979137
$ws.Range("I646").Value = 0.2782908341552996
$ws.Range("J646").Value = 0.5081153680062789
$ws.Range("K646").Value = 0.5971731448763251
$ws.Range("L646").Value = 0.9290442697558957

$ws.Range("H648").Value = 0.454892789255641
$ws.Range("I648").Value = 0.31036909929810197
$ws.Range("J648").Value = 0.828018210055704
$ws.Range("K648").Value = 0.8848484848484849
$ws.Range("L648").Value = 0.9021967526265521

$ws.Range("H649").Value = 0.4681372043233213
$ws.Range("I649").Value = 0.3155437486279438
$ws.Range("J649").Value = 0.8308291261911909
$ws.Range("K649").Value = 0.88268156424581
$ws.Range("L649").Value = 0.9042329055736449

$ws.Range("H651").Value = 0.5368380738132266
$ws.Range("I651").Value = 0.3556543493117752
$ws.Range("J651").Value = 0.733045202648451
$ws.Range("K651").Value = 0.8435114503816794
$ws.Range("L651").Value = 0.9149474503697937

$ws.Range("H652").Value = 0.5423956433481129
$ws.Range("I652").Value = 0.3572692590664213
$ws.Range("J652").Value = 0.7290896809931678
$ws.Range("K652").Value = 0.85
$ws.Range("L652").Value = 0.9149805447470817

$ws.Range("H654").Value = 0.453299574190362
$ws.Range("I654").Value = 0.29344521206585483
$ws.Range("J654").Value = 0.5280575399523143
$ws.Range("K654").Value = 0.642570281124498
$ws.Range("L654").Value = 0.9310485516115871

$ws.Range("H655").Value = 0.4344188358897779
$ws.Range("I655").Value = 0.2796335746330619
$ws.Range("J655").Value = 0.5081153680062789
$ws.Range("K655").Value = 0.6182495344506518
$ws.Range("L655").Value = 0.9329631914456097

$ws.Range("H657").Value = 0.46175600473454226
$ws.Range("I657").Value = 0.31327337239877606
$ws.Range("J657").Value = 0.828018210055704
$ws.Range("K657").Value = 0.8867924528301887
$ws.Range("L657").Value = 0.9149017363098645

$ws.Range("H658").Value = 0.46027963150212203
$ws.Range("I658").Value = 0.31237545163086405
$ws.Range("J658").Value = 0.8281075015977591
$ws.Range("K658").Value = 0.8881987577639752
$ws.Range("L658").Value = 0.9152510020996374

$ws.Range("H660").Value = 0.5459411712078307
$ws.Range("I660").Value = 0.3590445874609075
$ws.Range("J660").Value = 0.733045202648451
$ws.Range("K660").Value = 0.848605577689243
$ws.Range("L660").Value = 0.9273645368032628

$ws.Range("H661").Value = 0.5431879160883067
$ws.Range("I661").Value = 0.35749843186536034
$ws.Range("J661").Value = 0.7265112051562631
$ws.Range("K661").Value = 0.8443579766536965
$ws.Range("L661").Value = 0.9280575539568345

$ws.Range("H663").Value = 0.461540377220559
$ws.Range("I663").Value = 0.2962473185622385
$ws.Range("J663").Value = 0.5280575399523143
$ws.Range("K663").Value = 0.675
$ws.Range("L663").Value = 0.9366

$ws.Range("H664").Value = 0.461540377220559
$ws.Range("I664").Value = 0.2962473185622385
$ws.Range("J664").Value = 0.5280575399523143
$ws.Range("K664").Value = 0.627906976744186
$ws.Range("L664").Value = 0.941140653541709
